$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.478.27"
$ws.Range("E2").Value = "  -0.45%  "
$ws.Range("D3").Value = "2.284.70"
$ws.Range("E3").Value = "  -0.39%  "
$ws.Range("E4").Value = "  -0.20%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "310.27"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.65%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "103.16"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.06%  "
$ws.Range("E7").Value = "  -1.15%  "
$ws.Range("E8").Value = "  -0.23%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.602"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.26%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "38.68"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.69%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0896"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.39%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.15"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.60%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.107"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.62%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.969"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.24%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.26"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.27%  "
$ws.Range("D16").Value = "2.630.42"
$ws.Range("E16").Value = "  -0.43%  "
$ws.Range("D17").Value = "2.276.54"
$ws.Range("E17").Value = "  -0.76%  "
$ws.Range("D18").Value = "42.424.15"
$ws.Range("E18").Value = "  -0.39%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.28"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.14%  "
$ws.Range("E20").Value = "  -1.31%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.49"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.47%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "73.24"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.01%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "269.77"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.07%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.40"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -5.11%  "
$ws.Range("E25").Value = "  -3.30%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.07%  "
$ws.Range("E27").Value = "  -1.56%  "
$ws.Range("B28").Value = "Filecoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.95"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +12.91%  "
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.28"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.70%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "22.31"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.20%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "35.79"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -6.73%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "164.18"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.84%  "
$ws.Range("E33").Value = "  -4.13%  "
$ws.Range("E34").Value = "  -2.72%  "
$ws.Range("E35").Value = "  +0.58%  "
$ws.Range("E36").Value = "  -2.77%  "
$ws.Range("E37").Value = "  -3.55%  "
$ws.Range("E38").Value = "  -3.07%  "
$ws.Range("E39").Value = "  +0.71%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.61"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.57%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "112.78"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +18.95%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.55"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.82%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "70.15"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.14%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.00"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.36%  "
$ws.Range("E45").Value = "  -0.41%  "
$ws.Range("E46").Value = "  -2.31%  "
$ws.Range("D47").Value = "1.728.86"
$ws.Range("E47").Value = "  +9.08%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "110.05"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.19%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "77.20"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -5.39%  "
$ws.Range("E50").Value = "  -2.95%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.12"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.84%  "
